# Task 4 edit script
# 1) "I'm Jagadish Mali, and I'm eager ..." -> split into 3 runs, replacing
#    the name with " Ritika Juyal ".
# 2) "The map chart concludes ..." paragraph -> drop the w:proofErr markers,
#    merge into two runs, and move w:lastRenderedPageBreak to the second run.

$d = $word.ActiveDocument

function Insert-RunsXml($Range, $InnerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($pkg)
}

# --- Edit 1: introduce the new name ------------------------------------
$found = $d.Content
$found.Find.Execute("I'm Jagadish Mali, and I'm eager")
$para1 = $found.Paragraphs(1)
$range1 = $d.Range($para1.Range.Start, $para1.Range.End - 1)

$run1 = '<w:r><w:rPr><w:color w:val="252525"/></w:rPr><w:t>I''m</w:t></w:r>'
$run2 = '<w:r><w:rPr><w:color w:val="252525"/></w:rPr><w:t xml:space="preserve"> Ritika Juyal </w:t></w:r>'
$run3 = '<w:r><w:rPr><w:color w:val="252525"/></w:rPr><w:t>, and I''m eager to share some information with you regarding your company. I appreciate you giving me the leading questions. Seeing the sorts of insights, you expect to derive from the data was useful. I really believe you will find the analysis convincing and useful as you decide how to proceed with your next business prospects.</w:t></w:r>'

$combo1 = $run1 + $run2 + $run3
Insert-RunsXml $range1 $combo1

# --- Edit 2: rewrite the "map chart" paragraph --------------------------
$found2 = $d.Content
$found2.Find.Execute("The map chart concludes")
$para2 = $found2.Paragraphs(1)
$range2 = $d.Range($para2.Range.Start, $para2.Range.End - 1)

$run4 = '<w:r><w:rPr><w:color w:val="252525"/></w:rPr><w:t xml:space="preserve">The map chart concludes by comparing the places that have produced the greatest revenue to those that have not. Apart from the UK, it is clear that nations like the Netherlands, Ireland, Germany, France, and Australia generate large profits, and the company should invest more in </w:t></w:r>'
$run5 = '<w:r><w:rPr><w:color w:val="252525"/></w:rPr><w:lastRenderedPageBreak/><w:t>these nations to boost product demand. The map also reveals that the majority of sales occur only in the European zone, with only a small number in the American region. Along with Russia, there is no market for the items in Africa or Asia. Sales revenues and profitability might increase with the implementation of a fresh strategy focused on these areas.</w:t></w:r>'

$combo2 = $run4 + $run5
Insert-RunsXml $range2 $combo2
